$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fix row 2 Result column: SKIP -> PASS
$ws.Range("E2").Value = "PASS"

# Add new row 13: TestCase_F12
$ws.Range("A13").Value = "TestCase_F12"
$ws.Range("B13").Value = "OPQA-1183"
$ws.Range("C13").Value = "Verify that user is receiving notification when someone he is following created a public watch list. (single event notification)"
$ws.Range("D13").Value = "Y"
$ws.Range("E13").Value = "PASS"

# Add new row 14: TestCase_F13
$ws.Range("A14").Value = "TestCase_F13"
$ws.Range("C14").Value = "Verify that user is receiving notification when someone he is following made an existing watch list from private to public. (single event notification)"
$ws.Range("B14").Value = "OPQA-1184"
$ws.Range("D14").Value = "Y"
$ws.Range("E14").Value = "PASS"
